$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 1051.6666
$ws.Range("I11").Value = 1051.6666
$ws.Range("K11").Value = 1051.6666
$ws.Range("M11").Value = -911.6666
$ws.Range("H40").Value = 3368.7856
$ws.Range("I40").Value = 1318
$ws.Range("J40").Value = 4189.1
$ws.Range("K40").Value = 1318
$ws.Range("L40").Value = 4189.1
$ws.Range("M40").Value = -1143
$ws.Range("N40").Value = -4539.1
$ws.Range("H43").Value = 7330.8887
$ws.Range("J43").Value = 9496.333000000001
$ws.Range("L43").Value = 9496.333000000001
$ws.Range("N43").Value = -9634.333000000001
$ws.Range("H51").Value = 19583.166
$ws.Range("J51").Value = 4374.75
$ws.Range("L51").Value = 4374.75
$ws.Range("N51").Value = -5342.75
$ws.Range("H62").Value = 7455.6665
$ws.Range("J62").Value = 9086
$ws.Range("L62").Value = 9086
$ws.Range("N62").Value = -10334
$ws.Range("H64").Value = 8562.888999999999
$ws.Range("I64").Value = 4079.8
$ws.Range("J64").Value = 9581.772000000001
$ws.Range("K64").Value = 4079.8
$ws.Range("L64").Value = 9581.772000000001
$ws.Range("M64").Value = -3831.8
$ws.Range("N64").Value = -10077.772
$ws.Range("H65").Value = 7455.6665
$ws.Range("J65").Value = 9086
$ws.Range("L65").Value = 45430
$ws.Range("N65").Value = -51670
$ws.Range("H67").Value = 8562.888999999999
$ws.Range("I67").Value = 4079.8
$ws.Range("J67").Value = 9581.772000000001
$ws.Range("K67").Value = 4079.8
$ws.Range("L67").Value = 9581.772000000001
$ws.Range("M67").Value = -3221.8
$ws.Range("N67").Value = -11297.772
$ws.Range("H70").Value = 42862576
$ws.Range("I70").Value = 20003870
$ws.Range("J70").Value = 55561856
$ws.Range("K70").Value = 60011610
$ws.Range("L70").Value = 166685568
$ws.Range("M70").Value = -60011340
$ws.Range("N70").Value = -166686108
$ws.Range("H73").Value = 42862576
$ws.Range("I73").Value = 20003870
$ws.Range("J73").Value = 55561856
$ws.Range("K73").Value = 60011610
$ws.Range("L73").Value = 166685568
$ws.Range("M73").Value = -60010674
$ws.Range("N73").Value = -166687440
$ws.Range("H74").Value = 7458.6
$ws.Range("I74").Value = 6205.643
$ws.Range("K74").Value = 6205.643
$ws.Range("M74").Value = -5269.643
$ws.Range("H76").Value = 2998
$ws.Range("I76").Value = 2998
$ws.Range("K76").Value = 2998
$ws.Range("M76").Value = -2683
$ws.Range("H77").Value = 7458.6
$ws.Range("I77").Value = 6205.643
$ws.Range("K77").Value = 31028.215
$ws.Range("M77").Value = -26348.215
$ws.Range("H79").Value = 2998
$ws.Range("I79").Value = 2998
$ws.Range("K79").Value = 2998
$ws.Range("M79").Value = -1906
$ws.Range("H86").Value = 1891.25
$ws.Range("I86").Value = 1373.7142
$ws.Range("J86").Value = 2293.7778
$ws.Range("K86").Value = 1373.7142
$ws.Range("L86").Value = 2293.7778
$ws.Range("M86").Value = -250.7141999999999
$ws.Range("N86").Value = -4539.7778
$ws.Range("H89").Value = 1891.25
$ws.Range("I89").Value = 1373.7142
$ws.Range("J89").Value = 2293.7778
$ws.Range("K89").Value = 6868.571
$ws.Range("L89").Value = 11468.889
$ws.Range("M89").Value = -1252.571
$ws.Range("N89").Value = -22700.889
$ws.Range("H98").Value = 2016.16
$ws.Range("I98").Value = 1632.3182
$ws.Range("K98").Value = 1632.3182
$ws.Range("M98").Value = -134.3181999999999
$ws.Range("H100").Value = 3820.8965
$ws.Range("I100").Value = 2705.1538
$ws.Range("J100").Value = 4727.4375
$ws.Range("K100").Value = 2705.1538
$ws.Range("L100").Value = 4727.4375
$ws.Range("M100").Value = -2164.1538
$ws.Range("N100").Value = -5809.4375
$ws.Range("H103").Value = 1264.2222
$ws.Range("I103").Value = 1225.1538
$ws.Range("J103").Value = 1286.3043
$ws.Range("K103").Value = 3675.4614
$ws.Range("L103").Value = 3858.9129
$ws.Range("M103").Value = -3089.4614
$ws.Range("N103").Value = -5030.9129
$ws.Range("H122").Value = 2016.16
$ws.Range("I122").Value = 1632.3182
$ws.Range("K122").Value = 4896.9546
$ws.Range("M122").Value = -2446.9546
$ws.Range("H125").Value = 13928
$ws.Range("I125").Value = 2217.75
$ws.Range("J125").Value = 42032.6
$ws.Range("K125").Value = 19959.75
$ws.Range("L125").Value = 378293.4
$ws.Range("M125").Value = -17499.75
$ws.Range("N125").Value = -383213.4
$ws.Range("H129").Value = 1419.48
$ws.Range("J129").Value = 2152.4614
$ws.Range("L129").Value = 6457.3842
$ws.Range("N129").Value = -16457.3842
$ws.Range("H132").Value = 2422.4517
$ws.Range("I132").Value = 1780.8889
$ws.Range("K132").Value = 5342.6667
$ws.Range("M132").Value = -2812.6667
$ws.Range("H137").Value = 14567.154
$ws.Range("J137").Value = 5243
$ws.Range("L137").Value = 15729
$ws.Range("N137").Value = -20829
$ws.Range("H138").Value = 3578.475
$ws.Range("J138").Value = 4227.3887
$ws.Range("L138").Value = 12682.1661
$ws.Range("N138").Value = -22962.1661

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1463.8462
$ws.Range("I2").Value = 519.1667
$ws.Range("J2").Value = 3589.375
$ws.Range("K2").Value = 519.1667
$ws.Range("L2").Value = 3589.375
$ws.Range("M2").Value = -406.1667
$ws.Range("N2").Value = -3815.375
$ws.Range("H24").Value = 75000
$ws.Range("J24").Value = 75000
$ws.Range("L24").Value = 75000
$ws.Range("N24").Value = -75748
$ws.Range("H32").Value = 2902.0852
$ws.Range("I32").Value = 3029.5476
$ws.Range("J32").Value = 1831.4
$ws.Range("K32").Value = 3029.5476
$ws.Range("L32").Value = 1831.4
$ws.Range("M32").Value = -2742.5476
$ws.Range("N32").Value = -2405.4
$ws.Range("H61").Value = 4601.5
$ws.Range("I61").Value = 4701.5713
$ws.Range("K61").Value = 4701.5713
$ws.Range("M61").Value = -4489.5713
$ws.Range("H63").Value = 4876
$ws.Range("I63").Value = 2314
$ws.Range("K63").Value = 2314
$ws.Range("M63").Value = -1628
$ws.Range("H66").Value = 4876
$ws.Range("I66").Value = 2314
$ws.Range("K66").Value = 11570
$ws.Range("M66").Value = -8138
$ws.Range("H74").Value = 2136.44
$ws.Range("I74").Value = 1745.9546
$ws.Range("K74").Value = 1745.9546
$ws.Range("M74").Value = -871.9546
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").Value = ""
$ws.Range("H77").Value = 2136.44
$ws.Range("I77").Value = 1745.9546
$ws.Range("K77").Value = 8729.773000000001
$ws.Range("M77").Value = -4361.773000000001
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").Value = ""
$ws.Range("H97").Value = 443.36365
$ws.Range("I97").Value = 521.2941
$ws.Range("J97").Value = 178.4
$ws.Range("K97").Value = 521.2941
$ws.Range("L97").Value = 178.4
$ws.Range("M97").Value = -25.29409999999996
$ws.Range("N97").Value = -1170.4
$ws.Range("H100").Value = 75000
$ws.Range("J100").Value = 75000
$ws.Range("L100").Value = 75000
$ws.Range("N100").Value = -77164
$ws.Range("H102").Value = 4425.227
$ws.Range("I102").Value = 2404
$ws.Range("J102").Value = 8756.429
$ws.Range("K102").Value = 2404
$ws.Range("L102").Value = 8756.429
$ws.Range("M102").Value = -782
$ws.Range("N102").Value = -12000.429
$ws.Range("H106").Value = 47000
$ws.Range("J106").Value = 47000
$ws.Range("L106").Value = 47000
$ws.Range("N106").Value = -49524
$ws.Range("H110").Value = 352.6129
$ws.Range("I110").Value = 365.5862
$ws.Range("J110").Value = 164.5
$ws.Range("K110").Value = 365.5862
$ws.Range("L110").Value = 164.5
$ws.Range("M110").Value = 1679.4138
$ws.Range("N110").Value = -4254.5
$ws.Range("H116").Value = 1463.8462
$ws.Range("I116").Value = 519.1667
$ws.Range("J116").Value = 3589.375
$ws.Range("K116").Value = 519.1667
$ws.Range("L116").Value = 3589.375
$ws.Range("M116").Value = 1774.8333
$ws.Range("N116").Value = -8177.375
$ws.Range("H122").Value = 2992.5
$ws.Range("I122").Value = 2926.25
$ws.Range("K122").Value = 8778.75
$ws.Range("M122").Value = -6328.75
$ws.Range("H132").Value = 2594.111
$ws.Range("I132").Value = 2408.303
$ws.Range("K132").Value = 7224.909
$ws.Range("M132").Value = -4694.909
$ws.Range("H136").Value = 4601.5
$ws.Range("I136").Value = 4701.5713
$ws.Range("K136").Value = 14104.7139
$ws.Range("M136").Value = -11554.7139

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1463.8462
$ws.Range("I3").Value = 519.1667
$ws.Range("J3").Value = 3589.375
$ws.Range("K3").Value = 519.1667
$ws.Range("L3").Value = 3589.375
$ws.Range("M3").Value = -405.1667
$ws.Range("N3").Value = -3817.375
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").Value = ""
$ws.Range("H48").Value = 199999
$ws.Range("J48").Value = 199999
$ws.Range("L48").Value = 199999
$ws.Range("N48").Value = -200829
$ws.Range("H62").Value = 50000
$ws.Range("J62").Value = 50000
$ws.Range("L62").Value = 50000
$ws.Range("N62").Value = -51372
$ws.Range("H65").Value = 50000
$ws.Range("J65").Value = 50000
$ws.Range("L65").Value = 150000
$ws.Range("N65").Value = -156864
$ws.Range("H80").Value = 364.2381
$ws.Range("J80").Value = 397.1111
$ws.Range("L80").Value = 397.1111
$ws.Range("N80").Value = -2393.1111
$ws.Range("H83").Value = 364.2381
$ws.Range("J83").Value = 397.1111
$ws.Range("L83").Value = 1985.5555
$ws.Range("N83").Value = -11969.5555
$ws.Range("H86").Value = 2779.6667
$ws.Range("I86").Value = 2664.3
$ws.Range("J86").Value = 3356.5
$ws.Range("K86").Value = 2664.3
$ws.Range("L86").Value = 3356.5
$ws.Range("M86").Value = -1541.3
$ws.Range("N86").Value = -5602.5
$ws.Range("H89").Value = 2779.6667
$ws.Range("I89").Value = 2664.3
$ws.Range("J89").Value = 3356.5
$ws.Range("K89").Value = 13321.5
$ws.Range("L89").Value = 16782.5
$ws.Range("M89").Value = -7705.5
$ws.Range("N89").Value = -28014.5
$ws.Range("H107").Value = 5491
$ws.Range("I107").Value = 1397.25
$ws.Range("J107").Value = 8999.929
$ws.Range("K107").Value = 1397.25
$ws.Range("L107").Value = 8999.929
$ws.Range("M107").Value = 522.75
$ws.Range("N107").Value = -12839.929
$ws.Range("H134").Value = 4260.3125
$ws.Range("I134").Value = 4366.615
$ws.Range("K134").Value = 13099.845
$ws.Range("M134").Value = -10564.845

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 9750
$ws.Range("J4").Value = 10000
$ws.Range("L4").Value = 10000
$ws.Range("N4").Value = -10224
$ws.Range("H16").Value = 2747.25
$ws.Range("I16").Value = 1996.3334
$ws.Range("J16").Value = 5000
$ws.Range("K16").Value = 1996.3334
$ws.Range("L16").Value = 5000
$ws.Range("M16").Value = -1709.3334
$ws.Range("N16").Value = -5574
$ws.Range("H31").Value = 3082.8
$ws.Range("I31").Value = 1517.1364
$ws.Range("J31").Value = 5732.385
$ws.Range("K31").Value = 1517.1364
$ws.Range("L31").Value = 5732.385
$ws.Range("M31").Value = -1222.1364
$ws.Range("N31").Value = -6322.385
$ws.Range("H34").Value = 3082.8
$ws.Range("I34").Value = 1517.1364
$ws.Range("J34").Value = 5732.385
$ws.Range("K34").Value = 1517.1364
$ws.Range("L34").Value = 5732.385
$ws.Range("M34").Value = -1315.1364
$ws.Range("N34").Value = -6136.385
$ws.Range("H58").Value = 1723.4546
$ws.Range("I58").Value = 1506.5555
$ws.Range("K58").Value = 1506.5555
$ws.Range("M58").Value = -1303.5555
$ws.Range("H86").Value = 8552.727999999999
$ws.Range("I86").Value = 8156.0586
$ws.Range("J86").Value = 9901.4
$ws.Range("K86").Value = 8156.0586
$ws.Range("L86").Value = 9901.4
$ws.Range("M86").Value = -7033.0586
$ws.Range("N86").Value = -12147.4
$ws.Range("H89").Value = 8552.727999999999
$ws.Range("I89").Value = 8156.0586
$ws.Range("J89").Value = 9901.4
$ws.Range("K89").Value = 40780.29300000001
$ws.Range("L89").Value = 49507
$ws.Range("M89").Value = -35164.29300000001
$ws.Range("N89").Value = -60739
$ws.Range("H107").Value = 940.1177
$ws.Range("I107").Value = 955.1818
$ws.Range("J107").Value = 912.5
$ws.Range("K107").Value = 955.1818
$ws.Range("L107").Value = 912.5
$ws.Range("M107").Value = 964.8182
$ws.Range("N107").Value = -4752.5
$ws.Range("H113").Value = 2747.25
$ws.Range("I113").Value = 1996.3334
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 1996.3334
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = 173.6666
$ws.Range("N113").Value = -9340
$ws.Range("H132").Value = 2256.65
$ws.Range("I132").Value = 2256.65
$ws.Range("K132").Value = 6769.950000000001
$ws.Range("M132").Value = -4239.950000000001
$ws.Range("H134").Value = 1829.1034
$ws.Range("I134").Value = 1667.7778
$ws.Range("K134").Value = 5003.3334
$ws.Range("M134").Value = -2468.3334
$ws.Range("H136").Value = 1723.4546
$ws.Range("I136").Value = 1506.5555
$ws.Range("K136").Value = 4519.666499999999
$ws.Range("M136").Value = -1969.666499999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1053503.4
$ws.Range("I11").Value = 997.5333000000001
$ws.Range("J11").Value = 5000400
$ws.Range("K11").Value = 2992.5999
$ws.Range("L11").Value = 15001200
$ws.Range("M11").Value = -2852.5999
$ws.Range("N11").Value = -15001480
$ws.Range("H23").Value = 2816.4375
$ws.Range("I23").Value = 4217.2
$ws.Range("J23").Value = 2179.7273
$ws.Range("K23").Value = 12651.6
$ws.Range("L23").Value = 6539.1819
$ws.Range("M23").Value = -12416.6
$ws.Range("N23").Value = -7009.1819
$ws.Range("H34").Value = 322.44446
$ws.Range("J34").Value = 452.8
$ws.Range("L34").Value = 1358.4
$ws.Range("N34").Value = -1526.4
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").Value = ""
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").Value = ""
$ws.Range("H122").Value = 314.26923
$ws.Range("I122").Value = 140.66667
$ws.Range("J122").Value = 366.35
$ws.Range("K122").Value = 1266.00003
$ws.Range("L122").Value = 3297.15
$ws.Range("M122").Value = 1183.99997
$ws.Range("N122").Value = -8197.15
$ws.Range("H132").Value = 2453.6924
$ws.Range("I132").Value = 1839.6
$ws.Range("J132").Value = 2837.5
$ws.Range("K132").Value = 16556.4
$ws.Range("L132").Value = 25537.5
$ws.Range("M132").Value = -14026.4
$ws.Range("N132").Value = -30597.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 77.588234
$ws.Range("I2").Value = 46.076923
$ws.Range("K2").Value = 46.076923
$ws.Range("M2").Value = 66.92307700000001
$ws.Range("H49").Value = 27937.5
$ws.Range("J49").Value = 27937.5
$ws.Range("L49").Value = 27937.5
$ws.Range("N49").Value = -28305.5
$ws.Range("H70").Value = 9101.433000000001
$ws.Range("I70").Value = 7292.6816
$ws.Range("K70").Value = 7292.6816
$ws.Range("M70").Value = -7022.6816
$ws.Range("H73").Value = 9101.433000000001
$ws.Range("I73").Value = 7292.6816
$ws.Range("K73").Value = 7292.6816
$ws.Range("M73").Value = -6356.6816
$ws.Range("H80").Value = 11190.392
$ws.Range("I80").Value = 13836.846
$ws.Range("K80").Value = 13836.846
$ws.Range("M80").Value = -12838.846
$ws.Range("H83").Value = 11190.392
$ws.Range("I83").Value = 13836.846
$ws.Range("K83").Value = 69184.23
$ws.Range("M83").Value = -64192.23
$ws.Range("H97").Value = 5287.6
$ws.Range("I97").Value = 2036.24
$ws.Range("K97").Value = 2036.24
$ws.Range("M97").Value = -1540.24
$ws.Range("H102").Value = 4413.2856
$ws.Range("I102").Value = 4413.2856
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 4413.2856
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -2791.2856
$ws.Range("N102").Value = ""
$ws.Range("H113").Value = 32265036
$ws.Range("I113").Value = 71431290
$ws.Range("J113").Value = 10476.235
$ws.Range("K113").Value = 71431290
$ws.Range("L113").Value = 10476.235
$ws.Range("M113").Value = -71429120
$ws.Range("N113").Value = -14816.235
$ws.Range("H122").Value = 5235.359
$ws.Range("I122").Value = 5398.231
$ws.Range("J122").Value = 4909.615
$ws.Range("K122").Value = 16194.693
$ws.Range("L122").Value = 14728.845
$ws.Range("M122").Value = -13744.693
$ws.Range("N122").Value = -19628.845
$ws.Range("H132").Value = 3025.1936
$ws.Range("I132").Value = 2849.3076
$ws.Range("J132").Value = 3939.8
$ws.Range("K132").Value = 8547.9228
$ws.Range("L132").Value = 11819.4
$ws.Range("M132").Value = -6017.9228
$ws.Range("N132").Value = -16879.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 70000
$ws.Range("J36").Value = 70000
$ws.Range("L36").Value = 70000
$ws.Range("N36").Value = -71124
$ws.Range("H40").Value = 7459.1875
$ws.Range("I40").Value = 4836.75
$ws.Range("K40").Value = 4836.75
$ws.Range("M40").Value = -4700.75
$ws.Range("H42").Value = 13999
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 13999
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 13999
$ws.Range("M42").Value = ""
$ws.Range("N42").Value = -15125
$ws.Range("H46").Value = 2548.6206
$ws.Range("I46").Value = 1376.4546
$ws.Range("K46").Value = 1376.4546
$ws.Range("M46").Value = -1188.4546
$ws.Range("H49").Value = 13999
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 13999
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 13999
$ws.Range("M49").Value = ""
$ws.Range("N49").Value = -14293
$ws.Range("H55").Value = 362.7857
$ws.Range("I55").Value = 349
$ws.Range("J55").Value = 445.5
$ws.Range("K55").Value = 349
$ws.Range("L55").Value = 445.5
$ws.Range("M55").Value = -176
$ws.Range("N55").Value = -791.5
$ws.Range("H61").Value = 4681.05
$ws.Range("I61").Value = 2237.818
$ws.Range("K61").Value = 2237.818
$ws.Range("M61").Value = -2035.818
$ws.Range("H68").Value = 7130.3447
$ws.Range("I68").Value = 4864.6665
$ws.Range("J68").Value = 8149.9
$ws.Range("K68").Value = 4864.6665
$ws.Range("L68").Value = 8149.9
$ws.Range("M68").Value = -4115.6665
$ws.Range("N68").Value = -9647.9
$ws.Range("H71").Value = 7130.3447
$ws.Range("I71").Value = 4864.6665
$ws.Range("J71").Value = 8149.9
$ws.Range("K71").Value = 24323.3325
$ws.Range("L71").Value = 40749.5
$ws.Range("M71").Value = -20579.3325
$ws.Range("N71").Value = -48237.5
$ws.Range("H74").Value = 49084.855
$ws.Range("I74").Value = 47265.668
$ws.Range("K74").Value = 47265.668
$ws.Range("M74").Value = -46267.668
$ws.Range("H77").Value = 49084.855
$ws.Range("I77").Value = 47265.668
$ws.Range("K77").Value = 141797.004
$ws.Range("M77").Value = -136805.004
$ws.Range("H93").Value = 2369.9744
$ws.Range("I93").Value = 1683.9231
$ws.Range("K93").Value = 1683.9231
$ws.Range("M93").Value = -435.9231
$ws.Range("H113").Value = 4681.05
$ws.Range("I113").Value = 2237.818
$ws.Range("K113").Value = 2237.818
$ws.Range("M113").Value = -67.81800000000021
$ws.Range("H136").Value = 5027.3696
$ws.Range("I136").Value = 5110
$ws.Range("J136").Value = 4983.3
$ws.Range("K136").Value = 15330
$ws.Range("L136").Value = 14949.9
$ws.Range("M136").Value = -12780
$ws.Range("N136").Value = -20049.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1162.5
$ws.Range("I81").Value = 1150
$ws.Range("K81").Value = 2300
$ws.Range("M81").Value = -1239
$ws.Range("H84").Value = 1162.5
$ws.Range("I84").Value = 1150
$ws.Range("K84").Value = 11500
$ws.Range("M84").Value = -6196
$ws.Range("H96").Value = 4171.625
$ws.Range("J96").Value = 2949.5
$ws.Range("L96").Value = 2949.5
$ws.Range("N96").Value = -5695.5
$ws.Range("H107").Value = 1329.2632
$ws.Range("I107").Value = 1264.2222
$ws.Range("J107").Value = 2500
$ws.Range("K107").Value = 3792.6666
$ws.Range("L107").Value = 7500
$ws.Range("M107").Value = -1872.6666
$ws.Range("N107").Value = -11340
$ws.Range("H113").Value = 1039.2413
$ws.Range("I113").Value = 937
$ws.Range("K113").Value = 2811
$ws.Range("M113").Value = -641
$ws.Range("H136").Value = 1135.5428
$ws.Range("I136").Value = 1026.8966
$ws.Range("J136").Value = 1660.6666
$ws.Range("K136").Value = 3080.6898
$ws.Range("L136").Value = 4981.9998
$ws.Range("M136").Value = -530.6898000000001
$ws.Range("N136").Value = -10081.9998
